# Actualización automática del inventario: agrega el nuevo producto
# "Engranaje de cobre para cuchilla de Plotter HP" (fila 50) a Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "O0PI5Q"
$ws.Range("B50").Value = "Engranaje de cobre para cuchilla de Plotter HP"
$ws.Range("C50").Value = "T120 T125 T130 T210 T250 T230 T520 T525 T530 T630 T650 T730 T830"
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 250000
$ws.Range("F50").Value = 7
$ws.Range("G50").Value = 5
$ws.Range("H50").Formula = "=(E50-D50)*G50"
$ws.Range("I50").Formula = "=D50*F50"
$ws.Range("J50").Value = 0
